$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("华胜天成", "华胜天成", "*ST立方"),
    @("掌阅科技", "*ST立方", "华胜天成"),
    @("*ST立方", "深科技", "汉缆股份"),
    @("天奇股份", "贵州茅台", "巨力索具"),
    @("深科技", "东方财富", "协鑫集成"),
    @("汉缆股份", "捷成股份", "利欧股份"),
    @("五洲新春", "光线传媒", "博纳影业"),
    @("利欧股份", "掌阅科技", "掌阅科技"),
    @("光线传媒", "天奇股份", "光线传媒"),
    @("巨力索具", "利欧股份", "嘉美包装"),
    @("国安股份", "大位科技", "天奇股份"),
    @("大位科技", "汉缆股份", "克来机电"),
    @("协鑫集成", "巨力索具", "大位科技"),
    @("贵州茅台", "天汽模", "深科技"),
    @("东方财富", "亚星锚链", "紫金矿业"),
    @("捷成股份", "航发动力", "风语筑"),
    @("博纳影业", "二六三", "万向钱潮"),
    @("嘉美包装", "五洲新春", "航天发展"),
    @("天汽模", "协鑫集成", "百达精工"),
    @("锋龙股份", "国安股份", "豫能控股")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

